# "Primera Refactorizacion Cobranza por Cabeceras"
#
# - Toggle the "Usar" (C) column: rows 2-6 False->True, rows 10-14 True->False.
#   Done via copy/paste-special (values) from existing text cells so the
#   "True"/"False" strings stay plain text (shared-string), instead of
#   Value="True" which Excel auto-coerces to a native boolean.
# - Row 6 "database" value "facoep" -> "facoep1".
# - Apply a value-filter on the table's first column ("Parametros servidor")
#   for "password", which also hides the non-matching data rows.
# - Update the active selection (no more scrolled topLeftCell, new active cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Usar column: rows 2-6 -> True (copy the existing "True" text from C10) ---
foreach ($r in 2..6) {
    $ws.Range("C10").Copy()
    $ws.Cells.Item($r, 3).PasteSpecial(-4163)
}

# --- Usar column: rows 10-14 -> False (copy the existing "False" text from C7) ---
foreach ($r in 10..14) {
    $ws.Range("C7").Copy()
    $ws.Cells.Item($r, 3).PasteSpecial(-4163)
}

# --- database value for row 6 ---
$ws.Range("B6").Value = "facoep1"

# --- Filter the table by "password" in the first column ---
$lo = $ws.ListObjects.Item(1)
$lo.Range.AutoFilter(1, @("password"), 7)

# --- Update selection / scroll position ---
$ws.Range("A4").Select()
